$d = $word.ActiveDocument

# Update the date line
$null = $d.Content.Find.Execute("2026-01-27 Tuesday", $true, $false, $false, $false, $false, $true, 1, $false, "2026-01-28 Wednesday", 2)

# Update each arithmetic-problem cell in the table, addressed by row/column
# to avoid any accidental partial-text matches between similar expressions.
$t = $d.Tables(1)

$c = $t.Cell(1,1).Range
$c.End = $c.End - 1
$c.Text = "11+69="
$c = $t.Cell(1,2).Range
$c.End = $c.End - 1
$c.Text = "86-15="
$c = $t.Cell(1,3).Range
$c.End = $c.End - 1
$c.Text = "75-44="
$c = $t.Cell(1,4).Range
$c.End = $c.End - 1
$c.Text = "79-1="
$c = $t.Cell(1,5).Range
$c.End = $c.End - 1
$c.Text = "26+34="
$c = $t.Cell(2,1).Range
$c.End = $c.End - 1
$c.Text = "96-57="
$c = $t.Cell(2,2).Range
$c.End = $c.End - 1
$c.Text = "41+37="
$c = $t.Cell(2,3).Range
$c.End = $c.End - 1
$c.Text = "35+11="
$c = $t.Cell(2,4).Range
$c.End = $c.End - 1
$c.Text = "51-50="
$c = $t.Cell(2,5).Range
$c.End = $c.End - 1
$c.Text = "86+10="
$c = $t.Cell(3,1).Range
$c.End = $c.End - 1
$c.Text = "39-4="
$c = $t.Cell(3,2).Range
$c.End = $c.End - 1
$c.Text = "21-10="
$c = $t.Cell(3,3).Range
$c.End = $c.End - 1
$c.Text = "78-34="
$c = $t.Cell(3,4).Range
$c.End = $c.End - 1
$c.Text = "5+73="
$c = $t.Cell(3,5).Range
$c.End = $c.End - 1
$c.Text = "19+38="
$c = $t.Cell(4,1).Range
$c.End = $c.End - 1
$c.Text = "42+2="
$c = $t.Cell(4,2).Range
$c.End = $c.End - 1
$c.Text = "82-51="
$c = $t.Cell(4,3).Range
$c.End = $c.End - 1
$c.Text = "64-14="
$c = $t.Cell(4,4).Range
$c.End = $c.End - 1
$c.Text = "7+17="
$c = $t.Cell(4,5).Range
$c.End = $c.End - 1
$c.Text = "18+66="
$c = $t.Cell(5,1).Range
$c.End = $c.End - 1
$c.Text = "76-34="
$c = $t.Cell(5,2).Range
$c.End = $c.End - 1
$c.Text = "74-59="
$c = $t.Cell(5,3).Range
$c.End = $c.End - 1
$c.Text = "86+9="
$c = $t.Cell(5,4).Range
$c.End = $c.End - 1
$c.Text = "48+21="
$c = $t.Cell(5,5).Range
$c.End = $c.End - 1
$c.Text = "9+75="
$c = $t.Cell(6,1).Range
$c.End = $c.End - 1
$c.Text = "57-35="
$c = $t.Cell(6,2).Range
$c.End = $c.End - 1
$c.Text = "19-14="
$c = $t.Cell(6,3).Range
$c.End = $c.End - 1
$c.Text = "42+52="
$c = $t.Cell(6,4).Range
$c.End = $c.End - 1
$c.Text = "94-51="
$c = $t.Cell(6,5).Range
$c.End = $c.End - 1
$c.Text = "35+36="
$c = $t.Cell(7,1).Range
$c.End = $c.End - 1
$c.Text = "11-8="
$c = $t.Cell(7,2).Range
$c.End = $c.End - 1
$c.Text = "89-80="
$c = $t.Cell(7,3).Range
$c.End = $c.End - 1
$c.Text = "76-13="
$c = $t.Cell(7,4).Range
$c.End = $c.End - 1
$c.Text = "42+35="
$c = $t.Cell(7,5).Range
$c.End = $c.End - 1
$c.Text = "77-58="
$c = $t.Cell(8,1).Range
$c.End = $c.End - 1
$c.Text = "84-60="
$c = $t.Cell(8,2).Range
$c.End = $c.End - 1
$c.Text = "71-61="
$c = $t.Cell(8,3).Range
$c.End = $c.End - 1
$c.Text = "23-2="
$c = $t.Cell(8,4).Range
$c.End = $c.End - 1
$c.Text = "25+73="
$c = $t.Cell(8,5).Range
$c.End = $c.End - 1
$c.Text = "50-11="
$c = $t.Cell(9,1).Range
$c.End = $c.End - 1
$c.Text = "21+41="
$c = $t.Cell(9,2).Range
$c.End = $c.End - 1
$c.Text = "24+26="
$c = $t.Cell(9,3).Range
$c.End = $c.End - 1
$c.Text = "21+38="
$c = $t.Cell(9,4).Range
$c.End = $c.End - 1
$c.Text = "77-25="
$c = $t.Cell(9,5).Range
$c.End = $c.End - 1
$c.Text = "54+6="
$c = $t.Cell(10,1).Range
$c.End = $c.End - 1
$c.Text = "3+53="
$c = $t.Cell(10,2).Range
$c.End = $c.End - 1
$c.Text = "89-55="
$c = $t.Cell(10,3).Range
$c.End = $c.End - 1
$c.Text = "94-57="
$c = $t.Cell(10,4).Range
$c.End = $c.End - 1
$c.Text = "85-33="
$c = $t.Cell(10,5).Range
$c.End = $c.End - 1
$c.Text = "89-86="
$c = $t.Cell(11,1).Range
$c.End = $c.End - 1
$c.Text = "39-30="
$c = $t.Cell(11,2).Range
$c.End = $c.End - 1
$c.Text = "6+54="
$c = $t.Cell(11,3).Range
$c.End = $c.End - 1
$c.Text = "78-77="
$c = $t.Cell(11,4).Range
$c.End = $c.End - 1
$c.Text = "71-30="
$c = $t.Cell(11,5).Range
$c.End = $c.End - 1
$c.Text = "82-73="
$c = $t.Cell(12,1).Range
$c.End = $c.End - 1
$c.Text = "75-41="
$c = $t.Cell(12,2).Range
$c.End = $c.End - 1
$c.Text = "76-8="
$c = $t.Cell(12,3).Range
$c.End = $c.End - 1
$c.Text = "63+21="
$c = $t.Cell(12,4).Range
$c.End = $c.End - 1
$c.Text = "2+28="
$c = $t.Cell(12,5).Range
$c.End = $c.End - 1
$c.Text = "21-14="
$c = $t.Cell(13,1).Range
$c.End = $c.End - 1
$c.Text = "31+23="
$c = $t.Cell(13,2).Range
$c.End = $c.End - 1
$c.Text = "5+92="
$c = $t.Cell(13,3).Range
$c.End = $c.End - 1
$c.Text = "57-4="
$c = $t.Cell(13,4).Range
$c.End = $c.End - 1
$c.Text = "75-74="
$c = $t.Cell(13,5).Range
$c.End = $c.End - 1
$c.Text = "31+0="
$c = $t.Cell(14,1).Range
$c.End = $c.End - 1
$c.Text = "25+18="
$c = $t.Cell(14,2).Range
$c.End = $c.End - 1
$c.Text = "45+26="
$c = $t.Cell(14,3).Range
$c.End = $c.End - 1
$c.Text = "4+75="
$c = $t.Cell(14,4).Range
$c.End = $c.End - 1
$c.Text = "71-4="
$c = $t.Cell(14,5).Range
$c.End = $c.End - 1
$c.Text = "81+15="
$c = $t.Cell(15,1).Range
$c.End = $c.End - 1
$c.Text = "23+52="
$c = $t.Cell(15,2).Range
$c.End = $c.End - 1
$c.Text = "1+94="
$c = $t.Cell(15,3).Range
$c.End = $c.End - 1
$c.Text = "35+46="
$c = $t.Cell(15,4).Range
$c.End = $c.End - 1
$c.Text = "15+66="
$c = $t.Cell(15,5).Range
$c.End = $c.End - 1
$c.Text = "0+93="
$c = $t.Cell(16,1).Range
$c.End = $c.End - 1
$c.Text = "82-61="
$c = $t.Cell(16,2).Range
$c.End = $c.End - 1
$c.Text = "49+48="
$c = $t.Cell(16,3).Range
$c.End = $c.End - 1
$c.Text = "23-6="
$c = $t.Cell(16,4).Range
$c.End = $c.End - 1
$c.Text = "63-63="
$c = $t.Cell(16,5).Range
$c.End = $c.End - 1
$c.Text = "11+62="
$c = $t.Cell(17,1).Range
$c.End = $c.End - 1
$c.Text = "91-62="
$c = $t.Cell(17,2).Range
$c.End = $c.End - 1
$c.Text = "10+30="
$c = $t.Cell(17,3).Range
$c.End = $c.End - 1
$c.Text = "97-67="
$c = $t.Cell(17,4).Range
$c.End = $c.End - 1
$c.Text = "47+27="
$c = $t.Cell(17,5).Range
$c.End = $c.End - 1
$c.Text = "20+76="
$c = $t.Cell(18,1).Range
$c.End = $c.End - 1
$c.Text = "3+62="
$c = $t.Cell(18,2).Range
$c.End = $c.End - 1
$c.Text = "57+41="
$c = $t.Cell(18,3).Range
$c.End = $c.End - 1
$c.Text = "4+16="
$c = $t.Cell(18,4).Range
$c.End = $c.End - 1
$c.Text = "64-10="
$c = $t.Cell(18,5).Range
$c.End = $c.End - 1
$c.Text = "4+62="
$c = $t.Cell(19,1).Range
$c.End = $c.End - 1
$c.Text = "3+1="
$c = $t.Cell(19,2).Range
$c.End = $c.End - 1
$c.Text = "91-82="
$c = $t.Cell(19,3).Range
$c.End = $c.End - 1
$c.Text = "20+36="
$c = $t.Cell(19,4).Range
$c.End = $c.End - 1
$c.Text = "67-25="
$c = $t.Cell(19,5).Range
$c.End = $c.End - 1
$c.Text = "2+85="
$c = $t.Cell(20,1).Range
$c.End = $c.End - 1
$c.Text = "71+23="
$c = $t.Cell(20,2).Range
$c.End = $c.End - 1
$c.Text = "32+9="
$c = $t.Cell(20,3).Range
$c.End = $c.End - 1
$c.Text = "23+21="
$c = $t.Cell(20,4).Range
$c.End = $c.End - 1
$c.Text = "58-44="
$c = $t.Cell(20,5).Range
$c.End = $c.End - 1
$c.Text = "73-21="
